$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format first so values round-trip as strings,
# matching the original inlineStr cell typing (many look numeric/date-like).
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36", "E37", "D38", "E38", "B39", "C39", "D39", "E39", "B40", "C40", "D40", "E40", "E41", "B42", "C42", "D42", "E42", "B43", "C43", "D43", "E43", "B44", "C44", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "E51")
foreach ($ref in $targetCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "42.882.61"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.384.61"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "326.93"
$ws.Range("E5").Value = "  +5.76%  "
$ws.Range("D6").Value = "99.56"
$ws.Range("E6").Value = "  -8.14%  "
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "39.93"
$ws.Range("E10").Value = "  -8.94%  "
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "8.39"
$ws.Range("E12").Value = "  -4.45%  "
$ws.Range("D13").Value = "1.01"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "16.50"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "2.743.89"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "2.390.02"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "42.847.90"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "7.78"
$ws.Range("E19").Value = "  +7.62%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "3.75"
$ws.Range("E21").Value = "  +7.87%  "
$ws.Range("D22").Value = "75.36"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "273.03"
$ws.Range("E23").Value = "  +7.22%  "
$ws.Range("E24").Value = "  -7.29%  "
$ws.Range("D25").Value = "10.17"
$ws.Range("E25").Value = "  +13.31%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "11.49"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").Value = "23.87"
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "172.81"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "0.0905"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "35.27"
$ws.Range("E33").Value = "  -8.57%  "
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -8.80%  "
$ws.Range("E37").Value = "  -4.69%  "
$ws.Range("D38").Value = "3.84"
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  +7.88%  "
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").Value = "95.43"
$ws.Range("E42").Value = "  +50.82%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.227"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "68.85"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "116.53"
$ws.Range("E46").Value = "  +6.63%  "
$ws.Range("D47").Value = "11.87"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "5.43"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("D49").Value = "9.01"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "1.620.59"
$ws.Range("E50").Value = "  +9.53%  "
$ws.Range("E51").Value = "  -1.94%  "
